# Adding Master Data XLS
# Adds the French ("fra") translation rows for the app_detail master data
# table (rows 14-19), matching the existing id/name/descr/lang_code/
# is_active/cr_by/cr_dtimes layout used by the English ("eng") and
# Arabic ("ara") rows already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 14; Id = 10013; Name = "Pré-inscription";          Descr = "Portail Web pour les pré-inscriptions" },
    @{ Row = 15; Id = 10014; Name = "Client dinscription";      Descr = "Application de bureau pour les inscriptions" },
    @{ Row = 16; Id = 10015; Name = "Processeur dinscription";  Descr = "Demande de post-inscription" },
    @{ Row = 17; Id = 10016; Name = "Authentification ID";      Descr = "Application pour lauthentification du fournisseur de services tiers" },
    @{ Row = 18; Id = 10017; Name = "Contrôle didentité";       Descr = "Portail Web pour la configuration dapplications" },
    @{ Row = 19; Id = 10018; Name = "Portail Résident";         Descr = "Portail Web pour les services de génération de post-ID" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Id
    $ws.Range("B$row").Value = $r.Name
    $ws.Range("C$row").Value = $r.Descr
    $ws.Range("D$row").Value = "fra"
    $ws.Range("E$row").Value = $true
    $ws.Range("F$row").Value = "superadmin"
    $ws.Range("G$row").Value = "now()"
}

# Column sizing: column A best-fits the numeric id values, column B is
# widened to comfortably fit the longest name string.
$ws.Columns.Item(1).ColumnWidth = 5
$ws.Columns.Item(2).ColumnWidth = 18.5

# Scroll the view down and select the empty area below the table, as
# left by the editor after appending the rows.
$ws.Range("A20:XFD1048576").Select() | Out-Null

# Restore a "normal" (non-pagebreak-preview) print setup for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
